$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALXN")

# Row 6: "Change in inventories"
$ws.Range("B6").Value = -1800000.0
$ws.Range("C6").Value = 95000000.0
$ws.Range("D6").Value = 617700000.0
$ws.Range("E6").Value = 531000000.0
$ws.Range("F6").Value = 570600000.0
$ws.Range("G6").Value = 522100000.0

# Row 7: "Change in payables and accrued liability"
$ws.Range("B7").Value = 72300000.0
$ws.Range("C7").Value = 122500000.0
$ws.Range("D7").Value = 47300000.0
$ws.Range("E7").Value = 118600000.0
$ws.Range("F7").Value = 170700000.0
$ws.Range("G7").Value = 230700000.0

# Row 30: "Capital Stock Change" - B30 was empty inline string, now a number
$ws.Range("B30").Value = -391300000.0
